# Refresh the cryptocurrency "Price" (D) and "Volume(1h)" (E) columns with the
# latest scrape from the GitHub Actions job ("Updated symbol list").
#
# The source sheet stores these columns as literal text (e.g. "307.61",
# "0.27%") rather than numbers, so each new value is written back as text
# too. A plain `Range.Value = "..."` assignment lets Excel's smart-type
# detection reinterpret a numeric/percent-looking string as a Number (and
# silently round/trim it, e.g. "307.50" -> 307.5), and it also stamps a new
# number-format style onto the cell. To avoid both problems we:
#   1. Force the cell to a text number format ("@") before writing, so the
#      value is stored verbatim as text.
#   2. Immediately clear the format again so the cell's style index goes
#      back to the sheet's original (default) style - leaving only the
#      value changed, exactly like the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "307.50"
    "E2" = "0.60%"
    "D3" = "41.27"
    "E3" = "3.56%"
    "D4" = "5.127"
    "E4" = "1.70%"
    "D5" = "0.07635"
    "E5" = "-0.41%"
    "D6" = "4.262"
    "E6" = "0.39%"
    "D7" = "1.637"
    "E7" = "2.68%"
    "D8" = "2.464"
    "E8" = "0.08%"
    "D9" = "0.9023"
    "E9" = "2.29%"
    "D10" = "0.1099"
    "E10" = "10.12%"
    "D11" = "0.1776"
    "E11" = "3.01%"
    "D12" = "0.09203"
    "E12" = "2.93%"
    "D13" = "0.04187"
    "E13" = "-5.58%"
    "D14" = "0.1050"
    "E14" = "-0.40%"
    "D15" = "0.001252"
    "E15" = "-2.39%"
    "D16" = "0.005860"
    "E16" = "0.52%"
    "D17" = "3.356"
    "E17" = "-0.07%"
    "E18" = "-1.96%"
    "D19" = "6.547"
    "E19" = "-6.82%"
    "D20" = "0.1360"
    "E20" = "0.57%"
    "D21" = "0.2813"
    "E21" = "-13.11%"
    "D22" = "0.04120"
    "E22" = "-2.04%"
    "D23" = "0.001222"
    "E23" = "2.16%"
    "D24" = "0.003998"
    "E24" = "-1.45%"
    "E25" = "6.43%"
    "D38" = "0.02392"
    "E38" = "2.45%"
    "D39" = "0.05185"
    "E39" = "1.00%"
    "D40" = "0.007758"
    "E40" = "-2.31%"
    "E41" = "-1.49%"
    "D42" = "0.006961"
    "E42" = "4.33%"
    "E43" = "-1.89%"
    "D44" = "0.007706"
    "E44" = "-9.25%"
    "D45" = "0.3066"
    "E45" = "1.12%"
    "D46" = "0.00006860"
    "E46" = "4.91%"
    "D47" = "0.00000000750"
    "E47" = "-0.14%"
    "D48" = "0.01054"
    "E48" = "209.53%"
    "D49" = "0.004202"
    "E49" = "-40.08%"
    "E50" = "-0.14%"
    "E51" = "-0.14%"
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
    $range.ClearFormats()
}
